$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price / volume(1h) figures for the symbol list refresh.
$updates = @(
    @{ Cell = "D2"; Value = "300.81" },
    @{ Cell = "E2"; Value = "-1.15%" },
    @{ Cell = "D3"; Value = "31.49" },
    @{ Cell = "E3"; Value = "-3.17%" },
    @{ Cell = "D4"; Value = "5.158" },
    @{ Cell = "E4"; Value = "-2.66%" },
    @{ Cell = "D5"; Value = "0.07374" },
    @{ Cell = "E5"; Value = "-1.68%" },
    @{ Cell = "D6"; Value = "1.846" },
    @{ Cell = "E6"; Value = "22.86%" },
    @{ Cell = "D7"; Value = "7.885" },
    @{ Cell = "E7"; Value = "1.04%" },
    @{ Cell = "D8"; Value = "3.769" },
    @{ Cell = "E8"; Value = "-0.82%" },
    @{ Cell = "D9"; Value = "0.9269" },
    @{ Cell = "D10"; Value = "0.1707" },
    @{ Cell = "E10"; Value = "0.32%" },
    @{ Cell = "D11"; Value = "0.07414" },
    @{ Cell = "E11"; Value = "-4.79%" },
    @{ Cell = "D12"; Value = "0.08135" },
    @{ Cell = "E12"; Value = "0.73%" },
    @{ Cell = "D13"; Value = "0.03043" },
    @{ Cell = "E13"; Value = "0.26%" },
    @{ Cell = "D14"; Value = "0.09939" },
    @{ Cell = "E14"; Value = "0.24%" },
    @{ Cell = "D15"; Value = "0.001503" },
    @{ Cell = "E15"; Value = "-1.82%" },
    @{ Cell = "D16"; Value = "0.006066" },
    @{ Cell = "E16"; Value = "-5.73%" },
    @{ Cell = "D17"; Value = "3.451" },
    @{ Cell = "E17"; Value = "-0.65%" },
    @{ Cell = "D18"; Value = "2.225" },
    @{ Cell = "E18"; Value = "-0.18%" },
    @{ Cell = "D19"; Value = "0.3259" },
    @{ Cell = "E19"; Value = "-2.08%" },
    @{ Cell = "D20"; Value = "0.1346" },
    @{ Cell = "E20"; Value = "0.79%" },
    @{ Cell = "D21"; Value = "4.636" },
    @{ Cell = "E21"; Value = "1.59%" },
    @{ Cell = "D22"; Value = "0.04646" },
    @{ Cell = "E22"; Value = "1.04%" },
    @{ Cell = "D23"; Value = "0.1584" },
    @{ Cell = "E23"; Value = "-2.33%" },
    @{ Cell = "D24"; Value = "0.001219" },
    @{ Cell = "E24"; Value = "0.06%" },
    @{ Cell = "D25"; Value = "0.004477" },
    @{ Cell = "E25"; Value = "0.93%" },
    @{ Cell = "D26"; Value = "0.0001299" },
    @{ Cell = "E26"; Value = "-7.05%" },
    @{ Cell = "E27"; Value = "7.49%" },
    @{ Cell = "D39"; Value = "0.01716" },
    @{ Cell = "E39"; Value = "-2.66%" },
    @{ Cell = "D40"; Value = "0.04519" },
    @{ Cell = "E40"; Value = "-0.85%" },
    @{ Cell = "E41"; Value = "-1.26%" },
    @{ Cell = "D42"; Value = "0.1346" },
    @{ Cell = "E42"; Value = "0.04%" },
    @{ Cell = "D43"; Value = "0.002138" },
    @{ Cell = "E43"; Value = "-1.28%" },
    @{ Cell = "D44"; Value = "0.01043" },
    @{ Cell = "E44"; Value = "-17.87%" },
    @{ Cell = "D45"; Value = "0.00006269" },
    @{ Cell = "E45"; Value = "3.80%" },
    @{ Cell = "D46"; Value = "0.006994" },
    @{ Cell = "E46"; Value = "-46.21%" },
    @{ Cell = "D47"; Value = "0.7378" },
    @{ Cell = "E47"; Value = "4.03%" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    # Force text storage so numeric/percent-looking strings are not
    # reinterpreted as Number cells (matches the original inlineStr cells).
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Write-Host "Updated $($updates.Count) cells"
